$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the data rows 2-5 so that the former row 5 becomes the
# new row 2, and rows 2-4 each shift down by one row (row2->row3,
# row3->row4, row4->row5). Capture snapshots of every row first so that
# overwriting one row does not destroy data still needed for another.
$row2 = $ws.Range("A2:AY2").Value2()
$row3 = $ws.Range("A3:AY3").Value2()
$row4 = $ws.Range("A4:AY4").Value2()
$row5 = $ws.Range("A5:AY5").Value2()

# Columns Y and AA hold plain text dates (e.g. "2021-04-30") which Excel's
# automatic type detection would otherwise convert into date serial
# numbers when assigned through .Value2. Temporarily mark the destination
# columns as text so the strings survive the round trip unchanged, then
# restore the original "General" formatting afterwards.
$textCols = @("Y", "AA")
foreach ($col in $textCols) {
    $ws.Range(($col + "2") + ":" + ($col + "5")).NumberFormat = "@"
}

$ws.Range("A2:AY2").Value2 = $row5
$ws.Range("A3:AY3").Value2 = $row2
$ws.Range("A4:AY4").Value2 = $row3
$ws.Range("A5:AY5").Value2 = $row4

foreach ($col in $textCols) {
    $ws.Range(($col + "2") + ":" + ($col + "5")).NumberFormat = "General"
}

Write-Output "Row rotation complete"
